$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 47, shifting existing rows 47-108 down to 48-109.
$ws.Rows.Item(47).Insert()

# Populate the newly inserted row 47 with its data.
$ws.Cells.Item(47, 1).Value = 10
$ws.Cells.Item(47, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(47, 3).Value = "La Araucanía"
$ws.Cells.Item(47, 4).Value = 45225
$ws.Cells.Item(47, 4).NumberFormat = $ws.Cells.Item(46, 4).NumberFormat
$ws.Cells.Item(47, 5).Value = 9
$ws.Cells.Item(47, 6).Value = 100112042
$ws.Cells.Item(47, 7).Value = "Locoto"
$ws.Cells.Item(47, 8).Value = "Sin especificar"
$ws.Cells.Item(47, 9).Value = "Primera"
$ws.Cells.Item(47, 10).Value = 50
$ws.Cells.Item(47, 11).Value = 3800
$ws.Cells.Item(47, 12).Value = 3800
$ws.Cells.Item(47, 13).Value = 3800
$ws.Cells.Item(47, 14).Value = "$/kilo"
$ws.Cells.Item(47, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(47, 16).Value = 3800
$ws.Cells.Item(47, 17).Value = 1
$ws.Cells.Item(47, 18).Value = "Hortaliza"
